# Apply rank/data corrections to the expo location details sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 21
$ws.Range("A8").Value = 47
$ws.Range("A10").Value = 19
$ws.Range("E12").Value = 473

$ws.Range("A20").Value = 52
$ws.Range("C20").Value = "Fanfare Tickets"
$ws.Range("D20").Value = 4.4
$ws.Range("E20").Value = 14

$ws.Range("A21").Value = 30
$ws.Range("C21").Value = "Fern"
$ws.Range("D21").Value = 4.3
$ws.Range("E21").Value = 34

$ws.Range("A22").Value = 50
$ws.Range("C22").Value = "Floor & Decor"
$ws.Range("D22").Value = 4.4
$ws.Range("E22").Value = 281

$ws.Range("A23").Value = 37
$ws.Range("C23").Value = "Food Specialties Inc"
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0

$ws.Range("A25").Value = 24
$ws.Range("E29").Value = 528
$ws.Range("A33").Value = 22
$ws.Range("A38").Value = 23
$ws.Range("A41").Value = 25
$ws.Range("A47").Value = 51
$ws.Range("E49").Value = 1182
$ws.Range("E50").Value = 1757
$ws.Range("A53").Value = 20
